$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "ELIZABETH ACUÑA ALDANA" mora row (old row 17) ---
$ws.Rows.Item(17).Delete()

# --- Remove 5 of the 6 duplicated "MARIA JOSE OSPINO SIERRA" rows ---
# (old rows 20-24, i.e. rows 19-23 after the shift above; row 18 keeps the data)
$ws.Range("19:23").Delete()

# --- Remove the duplicate "QUELIS JOHANA..." row, keeping only the final
#     (bordered) "JESUS ALBERTO..." style row so the new 2-row block keeps
#     its original bottom border style ---
$ws.Rows.Item(19).Delete()

# --- Update header summary figures ---
$ws.Range("E11").Value = 180241
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 3

# --- Replace the remaining two data rows with the new account-holders ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1120740842"
$ws.Range("D18").Value = "CARLOS ALFREDO VIANA MONTEROSA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "19874875"
$ws.Range("D19").Value = "NELSON ARTURO TORRECILLA MOLINA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Columns.Item(4).ColumnWidth = 35.6328125
